# Generate Report for Handoff
# Adds two new handed-off files (45adacbf-760e-419c-8f27-a66d66377ffe.md and
# bb753406-7a5e-45b3-8cee-8cb79843d7a1.md) as rows 4 & 5 on the Overview,
# zh-cn and de-de worksheets/tables.

$wb = $excel.ActiveWorkbook

$commitHash = "968c9ee87147a502bd5325e0bafab96b7b74f6b1"
$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitHash/e2e/"

$file1 = "45adacbf-760e-419c-8f27-a66d66377ffe.md"
$file2 = "bb753406-7a5e-45b3-8cee-8cb79843d7a1.md"

$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E5").ColumnWidth = 17.2159881591797
$wsOverview.Range("F5").ColumnWidth = 17.2159881591797

# Row 4
$wsOverview.Range("A4").Value = $file1
$wsOverview.Range("B4").Value = "e2e\" + $file1
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $baseUrl + $file1, "", "", "e2e\" + $file1)
$wsOverview.Range("B4").Style = "HyperLink"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2017-02-17 07:58:49"
$wsOverview.Range("G4").NumberFormat = $dateFormat

# Row 5
$wsOverview.Range("A5").Value = $file2
$wsOverview.Range("B5").Value = "e2e\" + $file2
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), $baseUrl + $file2, "", "", "e2e\" + $file2)
$wsOverview.Range("B5").Style = "HyperLink"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2017-02-17 07:58:49"
$wsOverview.Range("G5").NumberFormat = $dateFormat

$loOverview = $wsOverview.ListObjects.Item("Overview")
$loOverview.Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C5").ColumnWidth = 17.2159881591797

# Row 4
$wsZhCn.Range("A4").Value = $file1
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), $baseUrl + $file1, "", "", $file1)
$wsZhCn.Range("A4").Style = "HyperLink"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "False"
$wsZhCn.Range("G4").Value = "45adacbf-760e-419c-8f27-a66d66377ffe.6e3d7dc9fc5e0000dbe65cfe718c97203585a820.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2017-02-17 07:58:33"
$wsZhCn.Range("H4").NumberFormat = $dateFormat
$wsZhCn.Range("I4").Value = ""
$wsZhCn.Range("J4").Value = ""
$wsZhCn.Range("K4").Value = ""
$wsZhCn.Range("L4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L4").NumberFormat = $dateFormat
$wsZhCn.Range("M4").Value = ""
$wsZhCn.Range("N4").Value = ""
$wsZhCn.Range("O4").Value = "True"
$wsZhCn.Range("P4").Value = ""
$wsZhCn.Range("Q4").Value = "False"
$wsZhCn.Range("R4").Value = ""

# Row 5
$wsZhCn.Range("A5").Value = $file2
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), $baseUrl + $file2, "", "", $file2)
$wsZhCn.Range("A5").Style = "HyperLink"
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "e2e"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("F5").Value = "False"
$wsZhCn.Range("G5").Value = "bb753406-7a5e-45b3-8cee-8cb79843d7a1.8c4aa768f8f2d60ed0c6cf02fd1b700b1ed8c121.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2017-02-17 07:58:33"
$wsZhCn.Range("H5").NumberFormat = $dateFormat
$wsZhCn.Range("I5").Value = ""
$wsZhCn.Range("J5").Value = ""
$wsZhCn.Range("K5").Value = ""
$wsZhCn.Range("L5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L5").NumberFormat = $dateFormat
$wsZhCn.Range("M5").Value = ""
$wsZhCn.Range("N5").Value = ""
$wsZhCn.Range("O5").Value = "True"
$wsZhCn.Range("P5").Value = ""
$wsZhCn.Range("Q5").Value = "False"
$wsZhCn.Range("R5").Value = ""

$loZhCn = $wsZhCn.ListObjects.Item("zh-cn")
$loZhCn.Resize($wsZhCn.Range("A1:R5"))

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C5").ColumnWidth = 17.2159881591797

# Row 4
$wsDeDe.Range("A4").Value = $file1
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), $baseUrl + $file1, "", "", $file1)
$wsDeDe.Range("A4").Style = "HyperLink"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "False"
$wsDeDe.Range("G4").Value = "45adacbf-760e-419c-8f27-a66d66377ffe.6e3d7dc9fc5e0000dbe65cfe718c97203585a820.de-de.xlf"
$wsDeDe.Range("H4").Value = "2017-02-17 07:58:49"
$wsDeDe.Range("H4").NumberFormat = $dateFormat
$wsDeDe.Range("I4").Value = ""
$wsDeDe.Range("J4").Value = ""
$wsDeDe.Range("K4").Value = ""
$wsDeDe.Range("L4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L4").NumberFormat = $dateFormat
$wsDeDe.Range("M4").Value = ""
$wsDeDe.Range("N4").Value = ""
$wsDeDe.Range("O4").Value = "True"
$wsDeDe.Range("P4").Value = ""
$wsDeDe.Range("Q4").Value = "False"
$wsDeDe.Range("R4").Value = ""

# Row 5
$wsDeDe.Range("A5").Value = $file2
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), $baseUrl + $file2, "", "", $file2)
$wsDeDe.Range("A5").Style = "HyperLink"
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "e2e"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("F5").Value = "False"
$wsDeDe.Range("G5").Value = "bb753406-7a5e-45b3-8cee-8cb79843d7a1.8c4aa768f8f2d60ed0c6cf02fd1b700b1ed8c121.de-de.xlf"
$wsDeDe.Range("H5").Value = "2017-02-17 07:58:49"
$wsDeDe.Range("H5").NumberFormat = $dateFormat
$wsDeDe.Range("I5").Value = ""
$wsDeDe.Range("J5").Value = ""
$wsDeDe.Range("K5").Value = ""
$wsDeDe.Range("L5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L5").NumberFormat = $dateFormat
$wsDeDe.Range("M5").Value = ""
$wsDeDe.Range("N5").Value = ""
$wsDeDe.Range("O5").Value = "True"
$wsDeDe.Range("P5").Value = ""
$wsDeDe.Range("Q5").Value = "False"
$wsDeDe.Range("R5").Value = ""

$loDeDe = $wsDeDe.ListObjects.Item("de-de")
$loDeDe.Resize($wsDeDe.Range("A1:R5"))
